$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B66 to be a numeric value (3) instead of text
$ws.Range("B66").Value = 3

# Add new row 67 with data
$ws.Range("A67").Value = "Ruilin"
$ws.Range("B67").Value = "'3"
$ws.Range("B67").Style = "Normal"
$ws.Range("C67").Value = "无"
$ws.Range("D67").Value = "SUG"
$ws.Range("E67").Value = "WRI"
$ws.Range("F67").Value = "b01bb119-e44b-4008-9381-38115d7c20f9"
$ws.Range("G67").Value = "mugzy2nI-Ayi1_annotated.xlsx"
$ws.Range("H67").Value = "It would be nice to have more explanation of the significance of beating SignalP."
